$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the abbreviation values between paired rows (same course/teacher),
# matching the reordering of shared strings in the diff:
#   row 9  (SYS)  <-> row 10 (KSYS)
#   row 12 (BIG)  <-> row 13 (KBIG)
#   row 17 (ANE)  <-> row 18 (KAEL)
$ws.Range("B9").Value  = "KSYS"
$ws.Range("B10").Value = "SYS"

$ws.Range("B12").Value = "KBIG"
$ws.Range("B13").Value = "BIG"

$ws.Range("B17").Value = "KAEL"
$ws.Range("B18").Value = "ANE"
